$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.160.25"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "3.908.41"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.74"
$ws.Range("E5").Value = "  +3.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.37"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -0.89%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("E9").Value = "  -0.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  -3.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000343"
$ws.Range("E11").Value = "  -4.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.92"
$ws.Range("E12").Value = "  -1.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.69"
$ws.Range("E13").Value = "  +2.79%  "

$ws.Range("D14").Value = "4.534.79"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("D15").Value = "3.903.07"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("E16").Value = "  -6.82%  "

$ws.Range("E17").Value = "  -1.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.80"
$ws.Range("E18").Value = "  -0.43%  "

$ws.Range("E19").Value = "  -2.21%  "

$ws.Range("D20").Value = "68.287.71"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.17"
$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.54"
$ws.Range("E22").Value = "  +4.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.85"
$ws.Range("E23").Value = "  +3.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.52"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.49"
$ws.Range("E25").Value = "  +17.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.25"
$ws.Range("E26").Value = "  +10.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.60"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.04"
$ws.Range("E28").Value = "  -1.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.72"
$ws.Range("E29").Value = "  +0.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "728.08"
$ws.Range("E30").Value = "  +1.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.71"
$ws.Range("E31").Value = "  +1.84%  "

$ws.Range("E32").Value = "  -2.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.24"
$ws.Range("E34").Value = "  +17.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.58"
$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("D36").Value = "0.0₃0864"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "60.33"
$ws.Range("E37").Value = "  +4.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.404"
$ws.Range("E38").Value = "  +20.25%  "

$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  +17.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0479"
$ws.Range("E42").Value = "  +0.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.16"
$ws.Range("E43").Value = "  +3.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.91"
$ws.Range("E44").Value = "  +3.33%  "

$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  +4.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.40"
$ws.Range("E48").Value = "  -3.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.12"
$ws.Range("E49").Value = "  -3.63%  "

$ws.Range("D50").Value = "0.0₆0343"
$ws.Range("E50").Value = "  +34.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.45"
$ws.Range("E51").Value = "  -2.23%  "
